$wb = $excel.ActiveWorkbook
Write-Host "Workbooks count:" $excel.Workbooks.Count
try {
    $wb2 = $excel.Workbooks.Add()
    Write-Host "Added wb2, count now:" $excel.Workbooks.Count
    Write-Host "wb2 name:" $wb2.Name
} catch {
    Write-Host "error adding workbook:" $_
}
